$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column O ("email" / gigibrantt@gmail.com), shifting P:R left to O:Q
$ws.Range("O1").EntireColumn.Delete()

# Update the saved selection to match the target state
$ws.Range("O1:O1048576").Select()
